# Clean-up of input tables:
#  - rename the worksheet from the stale "updated" name to "Tabelle1"
#  - move the selection/active cell back to A2 (also drops the scrolled
#    "topLeftCell" state left over from the previous session)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Tabelle1"

$ws.Range("A2").Select() | Out-Null
